# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect the latest scraped counts.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# row -> new value for column F
$updates = @{
    2  = 677
    3  = 516
    6  = 51
    8  = 3240
    9  = 4223
    10 = 109
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
